# The presentation currently uses the "Integral" (Red Violet) theme for its
# slide master (ppt/theme/theme2.xml) and the "Office Theme" colours for its
# notes master (ppt/theme/theme1.xml). The target edit swaps the two themes:
# the deck's visible theme becomes the default "Office Theme" colour scheme.
#
# This engine only exposes a writable theme-colour accessor for the part
# that is actually wired to a slide (ppt/theme/theme2.xml, via the slide
# master); we drive that accessor (Slide.ThemeColorScheme) to push the
# "Office Theme" colour values into every one of the 12 colour slots,
# recreating the effect of the theme swap described by the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Theme colour scheme slot order (PowerPoint COM ThemeColorScheme):
#  1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#  8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
# RGB ints are 0x00BBGGRR (i.e. R + G*256 + B*65536), matching the
# PowerPoint COM RGB() encoding.
$cs.Item(1).RGB  = 0          # dk1      000000
$cs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      44546A
$cs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  FFC000
$cs.Item(9).RGB  = 12874308   # accent5  4472C4
$cs.Item(10).RGB = 4697456    # accent6  70AD47
$cs.Item(11).RGB = 12673797   # hlink    0563C1
$cs.Item(12).RGB = 7491477    # folHlink 954F72
